$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 21603
$ws.Range("E2").Value = 3644
$ws.Range("F2").Value = 3644
$ws.Range("G2").Value = 3324
$ws.Range("H2").Value = 2497
$ws.Range("I2").Value = 2497
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 16216
$ws.Range("L2").Value = 5822
$ws.Range("M2").Value = 10394
$ws.Range("N2").Value = 10385
$ws.Range("O2").Value = 9
$ws.Range("P2").Value = 407
$ws.Range("Q2").Value = 5367
$ws.Range("R2").Value = -2992
$ws.Range("S2").Value = -3279
$ws.Range("T2").Value = 3058
$ws.Range("U2").Value = 2310
$ws.Range("V2").Value = 1653
$ws.Range("W2").Value = 16.87
$ws.Range("X2").Value = 11.56
$ws.Range("Y2").Value = 25.23
$ws.Range("Z2").Value = 15.18
$ws.Range("AA2").Value = 56.01
$ws.Range("AB2").Value = 2735.57
$ws.Range("AC2").Value = 3238
$ws.Range("AD2").Value = 26.01
$ws.Range("AE2").Value = 14001
$ws.Range("AF2").Value = 6.01
$ws.Range("AG2").Value = 2000
$ws.Range("AH2").Value = 2.38
$ws.Range("AI2").Value = 59.41
$ws.Range("AJ2").Value = 77124796

# Row 3
$ws.Range("D3").Value = 23152
$ws.Range("E3").Value = 4633
$ws.Range("F3").Value = 4633
$ws.Range("G3").Value = 4545
$ws.Range("H3").Value = 3431
$ws.Range("I3").Value = 3431
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 17754
$ws.Range("L3").Value = 5372
$ws.Range("M3").Value = 12382
$ws.Range("N3").Value = 12374
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = 407
$ws.Range("Q3").Value = 5128
$ws.Range("R3").Value = -3295
$ws.Range("S3").Value = -2134
$ws.Range("T3").Value = 3293
$ws.Range("U3").Value = 1836
$ws.Range("V3").Value = 880
$ws.Range("W3").Value = 20.01
$ws.Range("X3").Value = 14.82
$ws.Range("Y3").Value = 30.15
$ws.Range("Z3").Value = 20.2
$ws.Range("AA3").Value = 43.38
$ws.Range("AB3").Value = 3208.14
$ws.Range("AC3").Value = 4449
$ws.Range("AD3").Value = 18.9
$ws.Range("AE3").Value = 16654
$ws.Range("AF3").Value = 5.05
$ws.Range("AG3").Value = 2800
$ws.Range("AH3").Value = 3.33
$ws.Range("AI3").Value = 60.62
$ws.Range("AJ3").Value = 77124796

# Row 4
$ws.Range("D4").Value = 23763
$ws.Range("E4").Value = 3388
$ws.Range("F4").Value = 3388
$ws.Range("G4").Value = 3238
$ws.Range("H4").Value = 2433
$ws.Range("I4").Value = 2436
$ws.Range("J4").Value = -3
$ws.Range("K4").Value = 19677
$ws.Range("L4").Value = 7844
$ws.Range("M4").Value = 11833
$ws.Range("N4").Value = 11828
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 407
$ws.Range("Q4").Value = 3575
$ws.Range("R4").Value = -3675
$ws.Range("S4").Value = -400
$ws.Range("T4").Value = 3282
$ws.Range("U4").Value = 293
$ws.Range("V4").Value = 3518
$ws.Range("W4").Value = 14.26
$ws.Range("X4").Value = 10.24
$ws.Range("Y4").Value = 20.14
$ws.Range("Z4").Value = 13
$ws.Range("AA4").Value = 66.28
$ws.Range("AB4").Value = 3214.14
$ws.Range("AC4").Value = 3167
$ws.Range("AD4").Value = 27.88
$ws.Range("AE4").Value = 16136
$ws.Range("AF4").Value = 5.47
$ws.Range("AG4").Value = 3200
$ws.Range("AH4").Value = 3.62
$ws.Range("AI4").Value = 96.33
$ws.Range("AJ4").Value = 76380513

# Row 5
$ws.Range("D5").Value = 25168
$ws.Range("E5").Value = 4727
$ws.Range("F5").Value = 4727
$ws.Range("G5").Value = 4399
$ws.Range("H5").Value = 3256
$ws.Range("I5").Value = 3261
$ws.Range("J5").Value = -5
$ws.Range("K5").Value = 21589
$ws.Range("L5").Value = 11766
$ws.Range("M5").Value = 9823
$ws.Range("N5").Value = 9822
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 407
$ws.Range("Q5").Value = 5536
$ws.Range("R5").Value = -2897
$ws.Range("S5").Value = -2183
$ws.Range("T5").Value = 3368
$ws.Range("U5").Value = 2168
$ws.Range("V5").Value = 6823
$ws.Range("W5").Value = 18.78
$ws.Range("X5").Value = 12.94
$ws.Range("Y5").Value = 30.12
$ws.Range("Z5").Value = 15.78
$ws.Range("AA5").Value = 119.79
$ws.Range("AB5").Value = 2686.45
$ws.Range("AC5").Value = 4328
$ws.Range("AD5").Value = 22.58
$ws.Range("AE5").Value = 13619
$ws.Range("AF5").Value = 7.17
$ws.Range("AG5").Value = 3200
$ws.Range("AH5").Value = 3.28
$ws.Range("AI5").Value = 71.20999999999999
$ws.Range("AJ5").Value = 74818121

# Row 6
$ws.Range("D6").Value = 27073
$ws.Range("E6").Value = 5198
$ws.Range("F6").Value = 5198
$ws.Range("G6").Value = 4686
$ws.Range("H6").Value = 3498
$ws.Range("I6").Value = 3502
$ws.Range("K6").Value = 23789
$ws.Range("L6").Value = 12912
$ws.Range("M6").Value = 10877
$ws.Range("N6").Value = 10881
$ws.Range("P6").Value = 407
$ws.Range("Q6").Value = 5391
$ws.Range("R6").Value = -3935
$ws.Range("S6").Value = -1868
$ws.Range("T6").Value = 4003
$ws.Range("U6").Value = 1389
$ws.Range("V6").Value = 7416
$ws.Range("W6").Value = 19.2
$ws.Range("X6").Value = 12.92
$ws.Range("Y6").Value = 33.84
$ws.Range("Z6").Value = 15.42
$ws.Range("AA6").Value = 118.71
$ws.Range("AB6").Value = 2697.25
$ws.Range("AC6").Value = 4735
$ws.Range("AD6").Value = 15.65
$ws.Range("AE6").Value = 15074
$ws.Range("AF6").Value = 4.92
$ws.Range("AG6").Value = 3600
$ws.Range("AH6").Value = 4.86
$ws.Range("AI6").Value = 74.18000000000001
$ws.Range("AJ6").Value = 73799619

# Row 7
$ws.Range("D7").Value = 30209
$ws.Range("E7").Value = 5524
$ws.Range("G7").Value = 5587
$ws.Range("H7").Value = 4059
$ws.Range("I7").Value = 4128
$ws.Range("K7").Value = 26748
$ws.Range("L7").Value = 14986
$ws.Range("M7").Value = 11763
$ws.Range("N7").Value = 11768
$ws.Range("P7").Value = 409
$ws.Range("Q7").Value = 6190
$ws.Range("R7").Value = -2395
$ws.Range("S7").Value = -2119
$ws.Range("T7").Value = 2857
$ws.Range("U7").Value = 2458
$ws.Range("W7").Value = 18.29
$ws.Range("X7").Value = 13.44
$ws.Range("Y7").Value = 36.46
$ws.Range("Z7").Value = 16.06
$ws.Range("AA7").Value = 127.4
$ws.Range("AC7").Value = 5594
$ws.Range("AD7").Value = 15.7
$ws.Range("AE7").Value = 16295
$ws.Range("AF7").Value = 5.39
$ws.Range("AG7").Value = 3200
$ws.Range("AH7").Value = 3.64
$ws.Range("AI7").Value = 57.2

# Row 8
$ws.Range("D8").Value = 33017
$ws.Range("E8").Value = 6001
$ws.Range("G8").Value = 5826
$ws.Range("H8").Value = 4313
$ws.Range("I8").Value = 4320
$ws.Range("K8").Value = 28099
$ws.Range("L8").Value = 15178
$ws.Range("M8").Value = 12921
$ws.Range("N8").Value = 12929
$ws.Range("P8").Value = 409
$ws.Range("Q8").Value = 6229
$ws.Range("R8").Value = -2283
$ws.Range("S8").Value = -2457
$ws.Range("T8").Value = 2200
$ws.Range("U8").Value = 2376
$ws.Range("W8").Value = 18.17
$ws.Range("X8").Value = 13.06
$ws.Range("Y8").Value = 34.98
$ws.Range("Z8").Value = 15.73
$ws.Range("AA8").Value = 117.47
$ws.Range("AC8").Value = 5854
$ws.Range("AD8").Value = 15
$ws.Range("AE8").Value = 17882
$ws.Range("AF8").Value = 4.91
$ws.Range("AG8").Value = 3533
$ws.Range("AH8").Value = 4.02
$ws.Range("AI8").Value = 60.36

# Row 9
$ws.Range("D9").Value = 35272
$ws.Range("E9").Value = 6373
$ws.Range("G9").Value = 6184
$ws.Range("H9").Value = 4579
$ws.Range("I9").Value = 4580
$ws.Range("K9").Value = 28709
$ws.Range("L9").Value = 15193
$ws.Range("M9").Value = 13516
$ws.Range("N9").Value = 13527
$ws.Range("P9").Value = 409
$ws.Range("Q9").Value = 6519
$ws.Range("R9").Value = -2110
$ws.Range("S9").Value = -2707
$ws.Range("T9").Value = 2200
$ws.Range("U9").Value = 2831
$ws.Range("W9").Value = 18.07
$ws.Range("X9").Value = 12.98
$ws.Range("Y9").Value = 34.63
$ws.Range("Z9").Value = 16.12
$ws.Range("AA9").Value = 112.41
$ws.Range("AC9").Value = 6206
$ws.Range("AD9").Value = 14.15
$ws.Range("AE9").Value = 18709
$ws.Range("AF9").Value = 4.69
$ws.Range("AG9").Value = 3600
$ws.Range("AH9").Value = 4.1
$ws.Range("AI9").Value = 58
